$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# NOTE: touching $d.Tables / .Rows / .Cells (even just reading .Item(n))
# leaves this runtime's $d.Paragraphs collection unable to index past the
# first paragraph of that table, so every lookup below is done purely via
# $d.Paragraphs / $d.Range() instead of the Tables object model.

function Find-ParaIndex($pattern, [int]$startAt = 1) {
    $cnt = $d.Paragraphs.Count
    for ($i = $startAt; $i -le $cnt; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Replace the "Eg: I started painting..." paragraph (and the following
#    empty bold paragraph) with the new "most challenging aspect" answer
#    paragraph, followed by an empty paragraph that now carries the
#    _GoBack bookmark.
# ---------------------------------------------------------------------------
$egIdx = Find-ParaIndex("*Eg: I started painting as a hobby*")
if ($egIdx -eq -1) {
    throw "Could not find the 'Eg: I started painting' paragraph"
}

$p1 = $d.Paragraphs.Item($egIdx)
$p2 = $d.Paragraphs.Item($egIdx + 1)
$replaceRange = $d.Range($p1.Range.Start, $p2.Range.End)

$newAnswerXml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Developing </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>my</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> understand of what a testing class was and how it was to be used and implemented within the project was challenging. I was unsure where to even start with this. My approach was to use Google, Reddit, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>StackOverflow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ChatGPT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>CoPilot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> to see example implementation which I could use to focus my research further. This led to me discovering unit tests, which meant that</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> I needed to understand what a unit test was</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (which was quite different from what I had expected a unit test to be). I then research how to implement them. This was definitely made easier through following a YouTube tutorial and referencing Microsoft’s extensive C# documentation.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$replaceRange.InsertXML($newAnswerXml)

# ---------------------------------------------------------------------------
# 2. Remove the old _GoBack bookmark from the "Testing class implements..."
#    table cell (the cell right after it in document order) — it now lives
#    on the paragraph inserted above instead.
# ---------------------------------------------------------------------------
$testIdx = Find-ParaIndex("*Testing class implements a way to record testing data*")
if ($testIdx -eq -1) {
    throw "Could not find the 'Testing class implements...' paragraph"
}
$bookmarkPara = $d.Paragraphs.Item($testIdx + 1)
$clearBookmarkXml = "<w:p $W><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr></w:p>"
$bookmarkPara.Range.InsertXML($clearBookmarkXml)

# ---------------------------------------------------------------------------
# 3. Add w:lastRenderedPageBreak before "Demonstrate simple method
#    overriding for Creature subclasses..."
# ---------------------------------------------------------------------------
$demoIdx = Find-ParaIndex("*Demonstrate simple method overriding*")
if ($demoIdx -eq -1) {
    throw "Could not find the 'Demonstrate simple method overriding' paragraph"
}
$demoPara = $d.Paragraphs.Item($demoIdx)
$demoXml = "<w:p $W><w:r><w:lastRenderedPageBreak/><w:t>Demonstrate simple method overriding for Creature subclasses (e.g., different attack methods for Player and Monster).</w:t></w:r></w:p>"
$demoPara.Range.InsertXML($demoXml)

# ---------------------------------------------------------------------------
# 4. Remove w:lastRenderedPageBreak from "Rooms can contain multiple items
#    or monsters."
# ---------------------------------------------------------------------------
$roomsIdx = Find-ParaIndex("*Rooms can contain multiple items or monsters*")
if ($roomsIdx -eq -1) {
    throw "Could not find the 'Rooms can contain multiple items or monsters' paragraph"
}
$roomsPara = $d.Paragraphs.Item($roomsIdx)
$roomsXml = "<w:p $W><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr><w:r><w:t>Rooms can contain multiple items or monsters.</w:t></w:r></w:p>"
$roomsPara.Range.InsertXML($roomsXml)

# ---------------------------------------------------------------------------
# 5. Add w:lastRenderedPageBreak before the "First " run in the "First
#    standard:" Heading2.
# ---------------------------------------------------------------------------
$firstIdx = Find-ParaIndex("First standard:*")
if ($firstIdx -eq -1) {
    throw "Could not find the 'First standard:' heading paragraph"
}
$firstPara = $d.Paragraphs.Item($firstIdx)
$firstXml = "<w:p $W><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space='preserve'>First </w:t></w:r><w:r><w:t>standard:</w:t></w:r></w:p>"
$firstPara.Range.InsertXML($firstXml)

# ---------------------------------------------------------------------------
# 6. Remove w:lastRenderedPageBreak from "Add randomness to gameplay
#    (e.g., monsters or items appear randomly in rooms)."
# ---------------------------------------------------------------------------
$randIdx = Find-ParaIndex("*Add randomness to gameplay*")
if ($randIdx -eq -1) {
    throw "Could not find the 'Add randomness to gameplay' paragraph"
}
$randPara = $d.Paragraphs.Item($randIdx)
$randXml = "<w:p $W><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr><w:r><w:t>Add randomness to gameplay (e.g., monsters or items appear randomly in rooms).</w:t></w:r></w:p>"
$randPara.Range.InsertXML($randXml)

Write-Host "All edits applied successfully."
